$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1. Find the paragraph in the "Models which use this block:" table cell
#    that holds the BANSHEE_SYSTEM.mdl reference (and currently also
#    carries the stray "_GoBack" bookmark at its end).
# ----------------------------------------------------------------------
$targetIndex = -1
for ($p = 1; $p -le $d.Paragraphs.Count; $p++) {
    $txt = $d.Paragraphs.Item($p).Range.Text
    if ($txt -match "BANSHEE_SYSTEM\.mdl") {
        $targetIndex = $p
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not find the BANSHEE_SYSTEM.mdl paragraph"
}

$targetPara = $d.Paragraphs.Item($targetIndex)
$targetRange = $d.Range($targetPara.Range.Start, $targetPara.Range.End)

# Rewrite that paragraph without the bookmark, and immediately follow it
# with a brand-new paragraph (same run formatting) that references the
# new motor test model. Using InsertXML lets us reproduce the exact
# run/paragraph formatting (noProof + themed color) that a typed edit in
# Word would keep.
$newParagraphsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00332E86" w:rsidRDefault="00332E86" w:rsidP="00854A14"><w:pPr><w:rPr><w:noProof/><w:color w:val="1F497D" w:themeColor="text2"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:color w:val="1F497D" w:themeColor="text2"/></w:rPr><w:t>Energy/HIL/</w:t></w:r><w:r w:rsidR="00854A14"><w:rPr><w:noProof/><w:color w:val="1F497D" w:themeColor="text2"/></w:rPr><w:t>DistributionSystems/SimulinkOpal/Banshee/BANSHEE_SYSTEM.mdl</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:noProof/><w:color w:val="1F497D" w:themeColor="text2"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:color w:val="1F497D" w:themeColor="text2"/></w:rPr><w:t>Energy/HIL/Components/SimulinkOpal/Motor/IM200HP_Test_Model_v2.mdl</w:t></w:r></w:p>
'@
$targetRange.InsertXML($newParagraphsXml)

# ----------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark down into the empty paragraph that
#    follows the table (it was left stranded there after the last save).
# ----------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$bookmarkParagraphXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$lastRange.InsertXML($bookmarkParagraphXml)

Write-Output "Motor library documentation table updated."
